$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Ссылки и документация" cell that previously held the
# ГОС ШОО/toktom.kg/kao.kg reference text with the new text pointing
# editors to the MON KR website section.
$ws.Range("B23").Value = "На сайте МОН КР в разделе «Базисный учебный план. Предметные стандарты» Нужно указать сайт"

# Reflect the new active selection recorded in the saved workbook.
$ws.Range("B2").Select()
